# Managing Assets Using Orchestrator API
#
# - Renames sheet "Folders" -> "SpecificFolderAsset" and trims it down to a
#   single row (Name/Id pair for one specific asset).
# - Adds a new sheet "AllFolderAsset" right after it, listing the full set
#   of assets (header row + two data rows).

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (shared string), even when it looks
# like a number (e.g. "902087"), without leaving the cell's number format
# changed afterwards.
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

# 1) Rename the existing sheet and reduce it to the single asset row.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SpecificFolderAsset"

Set-TextValue $ws1.Cells.Item(1, 1) "TestAddAsset"
Set-TextValue $ws1.Cells.Item(1, 2) "902087"
$ws1.Range("A2:B2").ClearContents()

# 2) Add the new "AllFolderAsset" sheet right after "SpecificFolderAsset".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AllFolderAsset"

Set-TextValue $ws2.Cells.Item(1, 1) "Name"
Set-TextValue $ws2.Cells.Item(1, 2) "Id"
Set-TextValue $ws2.Cells.Item(2, 1) "TestAddAsset"
Set-TextValue $ws2.Cells.Item(2, 2) "902087"
Set-TextValue $ws2.Cells.Item(3, 1) "TestAddAssetCred"
Set-TextValue $ws2.Cells.Item(3, 2) "907838"
